$wb = $excel.ActiveWorkbook

# Update "Stand" dates on Geimpfte Personen sheet
$wsGP = $wb.Worksheets.Item("Geimpfte Personen")
$wsGP.Range("B1").Value = "Stand 20.1."
$wsGP.Range("C1").Value = "Stand 27.1."

# Update vaccination figures on Regional Geimpfte sheet
$wsRG = $wb.Worksheets.Item("Regional Geimpfte")

$wsRG.Range("B2").Value = "51,8"
$wsRG.Range("C2").Value = "15,8"
$wsRG.Range("D2").Value = " 7,2"
$wsRG.Range("E2").Value = "1,6"
$wsRG.Range("F2").Value = "1,6"
$wsRG.Range("G2").Value = "0,4"
$wsRG.Range("H2").Value = 1990889

$wsRG.Range("B3").Value = "33,0"
$wsRG.Range("C3").Value = " 7,8"
$wsRG.Range("D3").Value = "11,0"
$wsRG.Range("E3").Value = "2,4"
$wsRG.Range("F3").Value = "1,4"
$wsRG.Range("G3").Value = "0,3"
$wsRG.Range("H3").Value = 228138

$wsRG.Range("B4").Value = "41,2"
$wsRG.Range("C4").Value = "28,2"
$wsRG.Range("D4").Value = " 8,2"
$wsRG.Range("E4").Value = "3,9"
$wsRG.Range("F4").Value = "1,5"
$wsRG.Range("G4").Value = "0,7"
$wsRG.Range("H4").Value = 365003

$wsRG.Range("B5").Value = "50,6"
$wsRG.Range("C5").Value = "49,4"
$wsRG.Range("D5").Value = "20,5"
$wsRG.Range("E5").Value = "7,0"
$wsRG.Range("F5").Value = "1,6"
$wsRG.Range("G5").Value = "0,6"
$wsRG.Range("H5").Value = 103639

$wsRG.Range("B6").Value = "54,7"
$wsRG.Range("C6").Value = " 0,1"
$wsRG.Range("D6").Value = "12,8"
$wsRG.Range("E6").Value = "0,0"
$wsRG.Range("F6").Value = "2,5"
$wsRG.Range("G6").Value = "0,1"
$wsRG.Range("H6").Value = 68478

$wsRG.Range("B7").Value = "82,6"
$wsRG.Range("C7").Value = "17,4"
$wsRG.Range("D7").Value = " 6,6"
$wsRG.Range("E7").Value = "0,7"
$wsRG.Range("F7").Value = "2,1"
$wsRG.Range("G7").Value = "0,4"
$wsRG.Range("H7").Value = 19825

$wsRG.Range("B8").Value = "53,3"
$wsRG.Range("C8").Value = "12,4"
$wsRG.Range("D8").Value = "12,3"
$wsRG.Range("E8").Value = "1,6"
$wsRG.Range("F8").Value = "1,8"
$wsRG.Range("G8").Value = "0,3"
$wsRG.Range("H8").Value = 44231

$wsRG.Range("B9").Value = "38,9"
$wsRG.Range("C9").Value = "22,7"
$wsRG.Range("D9").Value = " 8,0"
$wsRG.Range("E9").Value = "1,9"
$wsRG.Range("F9").Value = "1,2"
$wsRG.Range("G9").Value = "0,5"
$wsRG.Range("H9").Value = 141478

$wsRG.Range("B10").Value = "82,2"
$wsRG.Range("C10").Value = "17,8"
$wsRG.Range("D10").Value = " 6,3"
$wsRG.Range("E10").Value = "0,6"
$wsRG.Range("F10").Value = "2,6"
$wsRG.Range("G10").Value = "0,6"
$wsRG.Range("H10").Value = 62212

$wsRG.Range("B11").Value = "55,0"
$wsRG.Range("C11").Value = " 8,4"
$wsRG.Range("D11").Value = " 4,1"
$wsRG.Range("E11").Value = "0,6"
$wsRG.Range("F11").Value = "1,5"
$wsRG.Range("G11").Value = "0,2"
$wsRG.Range("H11").Value = 152734

$wsRG.Range("B12").Value = "71,3"
$wsRG.Range("C12").Value = "19,0"
$wsRG.Range("D12").Value = " 0,0"
$wsRG.Range("E12").Value = "0,0"
$wsRG.Range("F12").Value = "1,3"
$wsRG.Range("G12").Value = "0,3"
$wsRG.Range("H12").Value = 357370

$wsRG.Range("B13").Value = "74,8"
$wsRG.Range("C13").Value = " 5,8"
$wsRG.Range("D13").Value = "18,4"
$wsRG.Range("E13").Value = "0,0"
$wsRG.Range("F13").Value = "3,2"
$wsRG.Range("G13").Value = "0,1"
$wsRG.Range("H13").Value = 139626

$wsRG.Range("B14").Value = "31,1"
$wsRG.Range("C14").Value = "17,5"
$wsRG.Range("D14").Value = "11,4"
$wsRG.Range("E14").Value = "5,5"
$wsRG.Range("F14").Value = "1,4"
$wsRG.Range("G14").Value = "0,5"
$wsRG.Range("H14").Value = 24255

$wsRG.Range("B15").Value = "37,7"
$wsRG.Range("C15").Value = " 4,0"
$wsRG.Range("D15").Value = " 2,1"
$wsRG.Range("E15").Value = "0,2"
$wsRG.Range("F15").Value = "1,7"
$wsRG.Range("G15").Value = "0,2"
$wsRG.Range("H15").Value = 86187

$wsRG.Range("B16").Value = "37,4"
$wsRG.Range("C16").Value = "22,7"
$wsRG.Range("D16").Value = " 5,5"
$wsRG.Range("E16").Value = "2,6"
$wsRG.Range("F16").Value = "1,4"
$wsRG.Range("G16").Value = "0,7"
$wsRG.Range("H16").Value = 59750

$wsRG.Range("B17").Value = "60,0"
$wsRG.Range("C17").Value = " 1,5"
$wsRG.Range("D17").Value = "11,4"
$wsRG.Range("E17").Value = "0,9"
$wsRG.Range("F17").Value = "2,8"
$wsRG.Range("G17").Value = "0,2"
$wsRG.Range("H17").Value = 91473

$wsRG.Range("B18").Value = "22,1"
$wsRG.Range("C18").Value = " 0,7"
$wsRG.Range("D18").Value = " 8,6"
$wsRG.Range("E18").Value = "0,2"
$wsRG.Range("F18").Value = "2,0"
$wsRG.Range("G18").Value = "0,1"
$wsRG.Range("H18").Value = 46490

